# Weekly update of "Hortaliza, Macroferia Regional de Talca - Repollo" data.
# A new week's record is added (row 161), which shifts the "Segunda" quality
# record that used to live in row 160 down to row 161, while the "Primera"
# quality values for the preceding weeks (rows 157-160) are refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 157: date + min/max/weighted/per-kg price refresh ---
$ws.Range("D157").Value = 44448
$ws.Range("K157").Value = 500
$ws.Range("L157").Value = 500
$ws.Range("M157").Value = 500
$ws.Range("P157").Value = 500

# --- Row 158: date, volume, min/max/weighted/per-kg price refresh ---
$ws.Range("D158").Value = 44167
$ws.Range("J158").Value = 3000
$ws.Range("K158").Value = 700
$ws.Range("L158").Value = 700
$ws.Range("M158").Value = 700
$ws.Range("P158").Value = 700

# --- Row 159: date + volume refresh ---
$ws.Range("D159").Value = 44399
$ws.Range("J159").Value = 5000

# --- Row 160: now holds the "Primera" record that used to be duplicated ---
$ws.Range("I160").Value = "Primera"
$ws.Range("J160").Value = 4000
$ws.Range("K160").Value = 400
$ws.Range("L160").Value = 400
$ws.Range("M160").Value = 400
$ws.Range("P160").Value = 400

# --- New row 161: the "Segunda" record that used to live in row 160 ---
$ws.Range("A161").Value = 5
$ws.Range("B161").Value = "Macroferia Regional de Talca"
$ws.Range("C161").Value = "Maule"
$ws.Range("D161").NumberFormat = $ws.Range("D160").NumberFormat
$ws.Range("D161").Value = 44400
$ws.Range("E161").Value = 7
$ws.Range("F161").Value = 100112006
$ws.Range("G161").Value = "Repollo"
$ws.Range("H161").Value = "Crespo record"
$ws.Range("I161").Value = "Segunda"
$ws.Range("J161").Value = 2000
$ws.Range("K161").Value = 300
$ws.Range("L161").Value = 300
$ws.Range("M161").Value = 300
$ws.Range("N161").Value = "$/unidad"
$ws.Range("O161").Value = "Región del Maule"
$ws.Range("P161").Value = 300
$ws.Range("Q161").Value = 1
$ws.Range("R161").Value = "Hortaliza"
